$d = $word.ActiveDocument

# --- paragraph index 0 (Paragraphs.Item(1)) ---
$p0 = $d.Paragraphs.Item(1)
$r0 = $p0.Range
$xml0 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="116546DF" w14:textId="5D249EEB" w:rsidR="00577161" w:rsidRPr="00502BD6" w:rsidRDefault="00577161" w:rsidP="00577161"><w:pPr><w:widowControl w:val="0"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:instrText xml:space="preserve"> SEQ CHAPTER \h \r 1</w:instrText></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">IN THE </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>MAGISTRATE</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> COURT </w:t></w:r><w:r w:rsidR="00434CB7" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">OF </w:t></w:r><w:r w:rsidR="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:r w:rsidR="000E0C2A"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>case_</w:t></w:r><w:r w:rsidR="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>county</w:t></w:r><w:r w:rsidR="000E0C2A"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>.upper()</w:t></w:r><w:r w:rsidR="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r><w:r w:rsidR="00434CB7" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>COUNTY, WEST VIRGINIA</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r0.InsertXML($xml0)

# --- paragraph index 2 (Paragraphs.Item(3)) ---
$p2 = $d.Paragraphs.Item(3)
$r2 = $p2.Range
$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="792D7899" w14:textId="784DA8A2" w:rsidR="00577161" w:rsidRPr="00502BD6" w:rsidRDefault="00502BD6" w:rsidP="00577161"><w:pPr><w:widowControl w:val="0"/><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="001600E1"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>p</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="001600E1"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>.name.full() }}</w:t></w:r><w:r w:rsidR="00577161" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r2.InsertXML($xml2)

# --- paragraph index 4 (Paragraphs.Item(5)) ---
$p4 = $d.Paragraphs.Item(5)
$r4 = $p4.Range
$xml4 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1CDBB6FF" w14:textId="352D0DC7" w:rsidR="00076ECF" w:rsidRPr="00502BD6" w:rsidRDefault="00577161" w:rsidP="00577161"><w:pPr><w:widowControl w:val="0"/><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>Plaintiff</w:t></w:r><w:r w:rsidR="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>/Appellee</w:t></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00FF4437"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Magistrate Court Case </w:t></w:r><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>No</w:t></w:r><w:r w:rsidR="00985550" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00076ECF"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00076ECF"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>{{ case</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00076ECF"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>_num }}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r4.InsertXML($xml4)

# --- paragraph index 8 (Paragraphs.Item(9)) ---
$p8 = $d.Paragraphs.Item(9)
$r8 = $p8.Range
$xml8 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="357974F1" w14:textId="023DB722" w:rsidR="00577161" w:rsidRPr="00502BD6" w:rsidRDefault="001600E1" w:rsidP="00577161"><w:pPr><w:widowControl w:val="0"/><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>{{ r</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>.name.full(middle=”full”) }}</w:t></w:r><w:r w:rsidR="00577161" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">,  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r8.InsertXML($xml8)

# --- paragraph index 18 (Paragraphs.Item(19)) ---
$p18 = $d.Paragraphs.Item(19)
$r18 = $p18.Range
$xml18 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="78661446" w14:textId="006154CC" w:rsidR="00C63A94" w:rsidRPr="00502BD6" w:rsidRDefault="00902AF3" w:rsidP="0030551E"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00502BD6"><w:rPr><w:b/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>COME</w:t></w:r><w:r w:rsidR="00A12106" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>S</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> NOW,</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the Defendant</w:t></w:r><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>/Appellant</w:t></w:r><w:r w:rsidR="00434CB7" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>{{ r</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>.name.full(middle=”full”) }}</w:t></w:r><w:r w:rsidR="00353E76" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:b/><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">having filed a Notice of Appeal in the Magistrate Court of </w:t></w:r><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:r w:rsidR="000E0C2A"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>case_</w:t></w:r><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>county }}</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> County, West Virginia, in Case No. </w:t></w:r><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>{{ case_num }}</w:t></w:r><w:r w:rsidR="004C1D9C" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">and </w:t></w:r><w:r w:rsidR="000710AD"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>provides notice to</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> this Court </w:t></w:r><w:r w:rsidR="00FC4F7D"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>and all parties of</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the automatic stay provided by West Virginia law. Such a</w:t></w:r><w:r w:rsidR="00FC4F7D"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> stay is</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to avoid the eviction of the Defendant</w:t></w:r><w:r w:rsidR="00D24A47"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and any other enforcement </w:t></w:r><w:r w:rsidR="00AB4B88"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>of the judgment</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> until </w:t></w:r><w:r w:rsidR="001600E1"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>their</w:t></w:r><w:r w:rsidR="00C63A94" w:rsidRPr="00502BD6"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> appeal has been decided.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r18.InsertXML($xml18)

# --- paragraph index 21 (Paragraphs.Item(22)) ---
$p21 = $d.Paragraphs.Item(22)
$r21 = $p21.Range
$xml21 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1C19F165" w14:textId="77777777" w:rsidR="00E4328C" w:rsidRDefault="00353E76" w:rsidP="00353E76"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">The Defendant has filed a timely appeal from the judgment entered </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00E4328C" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>{{ judgment</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00E4328C" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>_date }}</w:t></w:r><w:r w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r21.InsertXML($xml21)

# --- paragraph index 23 (Paragraphs.Item(24)) ---
$p23 = $d.Paragraphs.Item(24)
$r23 = $p23.Range
$xml23 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4FA6CBCB" w14:textId="612820DA" w:rsidR="00353E76" w:rsidRDefault="00E4328C" w:rsidP="00353E76"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>A</w:t></w:r><w:r w:rsidR="00353E76" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>ny eviction</w:t></w:r><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> attempt</w:t></w:r><w:r w:rsidR="00353E76" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> would be unlawful. The Defendant is entitled to possession of the rental property </w:t></w:r><w:r w:rsidR="00065A80" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">at </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>{{ r</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>.address.on_one_line() }}</w:t></w:r><w:r w:rsidR="00B721E2" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00353E76" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">until </w:t></w:r><w:r w:rsidR="009D5F7E"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>their</w:t></w:r><w:r w:rsidR="00353E76" w:rsidRPr="00E4328C"><w:rPr><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> appeal has been decided. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r23.InsertXML($xml23)

# --- paragraph index 28 (Paragraphs.Item(29)) ---
$p28 = $d.Paragraphs.Item(29)
$r28 = $p28.Range
$xml28 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2D0AF33C" w14:textId="622CB97B" w:rsidR="00E4328C" w:rsidRPr="00502BD6" w:rsidRDefault="00E4328C" w:rsidP="00E4328C"><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="480" w:lineRule="auto"/><w:jc w:val="right"/><w:rPr><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>{{ r</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:szCs w:val="24"/></w:rPr><w:t>.name.full(middle=”full”) }}, Defendant/Appellant</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r28.InsertXML($xml28)
